$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '26.087.46'
Set-TextValue 'E2' '  +0.54%  '
Set-TextValue 'D3' '1.640.49'
Set-TextValue 'E3' '  +0.10%  '
Set-TextValue 'E4' '  +0.73%  '
Set-TextValue 'D5' '214.45'
Set-TextValue 'E5' '  -0.54%  '
Set-TextValue 'D6' '0.505'
Set-TextValue 'E6' '  -0.08%  '
Set-TextValue 'E7' '  +0.75%  '
Set-TextValue 'D8' '0.250'
Set-TextValue 'E8' '  -2.54%  '
Set-TextValue 'D9' '0.0624'
Set-TextValue 'D10' '18.43'
Set-TextValue 'E10' '  -6.23%  '
Set-TextValue 'E11' '  -0.14%  '
Set-TextValue 'D12' '1.773.59'
Set-TextValue 'E12' '  +8.33%  '
Set-TextValue 'D13' '4.20'
Set-TextValue 'E13' '  -1.66%  '
Set-TextValue 'D14' '0.528'
Set-TextValue 'E14' '  -2.97%  '
Set-TextValue 'B15' 'WrappedBTC'
Set-TextValue 'C15' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D15' '26.091.61'
Set-TextValue 'E15' '  +0.29%  '
Set-TextValue 'B16' 'Litecoin'
Set-TextValue 'C16' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D16' '62.20'
Set-TextValue 'E16' '  -1.24%  '
Set-TextValue 'D17' '0.0₃0747'
Set-TextValue 'E17' '  -2.08%  '
Set-TextValue 'E18' '  +0.76%  '
Set-TextValue 'D19' '189.70'
Set-TextValue 'E19' '  -1.71%  '
Set-TextValue 'D20' '4.26'
Set-TextValue 'E20' '  -2.33%  '
Set-TextValue 'D21' '9.51'
Set-TextValue 'E21' '  -4.19%  '
Set-TextValue 'D22' '6.11'
Set-TextValue 'E22' '  -3.00%  '
Set-TextValue 'D23' '144.06'
Set-TextValue 'E23' '  +0.28%  '
Set-TextValue 'E24' '  +0.13%  '
Set-TextValue 'E25' '  +0.79%  '
Set-TextValue 'E26' '  -1.46%  '
Set-TextValue 'E27' '  -2.05%  '
Set-TextValue 'D28' '15.19'
Set-TextValue 'E28' '  -2.52%  '
Set-TextValue 'E29' '  -0.59%  '
Set-TextValue 'D30' '0.0483'
Set-TextValue 'E30' '  -3.91%  '
Set-TextValue 'D31' '3.17'
Set-TextValue 'E31' '  -2.39%  '
Set-TextValue 'D32' '3.16'
Set-TextValue 'E32' '  -4.16%  '
Set-TextValue 'E33' '  -0.07%  '
Set-TextValue 'E34' '  -2.06%  '
Set-TextValue 'D35' '0.879'
Set-TextValue 'E35' '  -2.46%  '
Set-TextValue 'D36' '1.122.77'
Set-TextValue 'E36' '  -1.20%  '
Set-TextValue 'E37' '  +0.29%  '
Set-TextValue 'D38' '0.520'
Set-TextValue 'E38' '  -4.61%  '
Set-TextValue 'E39' '  -1.75%  '
Set-TextValue 'D40' '98.78'
Set-TextValue 'E40' '  -0.40%  '
Set-TextValue 'D41' '0.786'
Set-TextValue 'E41' '  -1.50%  '
Set-TextValue 'D42' '5.29'
Set-TextValue 'E42' '  -3.87%  '
Set-TextValue 'E43' '  -0.39%  '
Set-TextValue 'D44' '55.03'
Set-TextValue 'E44' '  -2.95%  '
Set-TextValue 'E45' '  -2.00%  '
Set-TextValue 'E46' '  -0.29%  '
Set-TextValue 'E47' '  +0.29%  '
Set-TextValue 'E48' '  -1.43%  '
Set-TextValue 'E49' '  +0.62%  '
Set-TextValue 'D50' '0.0926'
Set-TextValue 'E50' '  -3.95%  '
Set-TextValue 'E51' '  -1.20%  '
